$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny precision change on existing row 49 (B49)
$ws.Range("B49").Value = 44708.07832020833

# New rows 50-58: Name (shared string) + Date (serial number)
$newRows = @(
    @("Tulsi",  44708.70932243056),
    @("Tulsi",  44708.70997974537),
    @("Ajay",   44708.73424809028),
    @("Ajay",   44708.73463961806),
    @("Aditya", 44708.74323366898),
    @("Ajay",   44709.00327535879),
    @("Aditya", 44709.00349340278),
    @("Ajay",   44709.00361327546),
    @("Ajay",   44709.01079778255)
)

$startRow = 50
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $name = $newRows[$i][0]
    $date = $newRows[$i][1]
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $date
    $ws.Cells.Item($row, 2).NumberFormat = $ws.Cells.Item($row - 1, 2).NumberFormat
}
